# Scheduled market-data refresh: update computed price/profit columns (H:N)
# across the crafting-class Leve sheets. Values sourced from the latest
# Universalis average-price snapshot; mirrors the commit's XML diff 1:1.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 388.6111
$ws.Range("I53").Value = 325.2
$ws.Range("J53").Value = 467.875
$ws.Range("K53").Value = 325.2
$ws.Range("L53").Value = 467.875
$ws.Range("M53").Value = 311.8
$ws.Range("N53").Value = -1741.875
# row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 794241.5600000001
$ws.Range("I107").Value = 1010694.44
$ws.Range("J107").Value = 581
$ws.Range("K107").Value = 1010694.44
$ws.Range("L107").Value = 581
$ws.Range("M107").Value = -1008774.44
$ws.Range("N107").Value = -4421
# row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 1097.9412
$ws.Range("I111").Value = 708
$ws.Range("J111").Value = 1812.8334
$ws.Range("K111").Value = 2124
$ws.Range("L111").Value = 5438.5002
$ws.Range("M111").Value = 943
$ws.Range("N111").Value = -11572.5002
# row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 166667970
$ws.Range("I137").Value = 200001170
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 600003510
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -600000960
$ws.Range("N137").Value = -11100
# row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 3140.125
$ws.Range("I141").Value = 2487.0417
$ws.Range("K141").Value = 7461.125100000001
$ws.Range("M141").Value = -2281.125100000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 40552.652
$ws.Range("I2").Value = 65125.312
$ws.Range("K2").Value = 65125.312
$ws.Range("M2").Value = -65012.312
# row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1683.5264
$ws.Range("I110").Value = 763.375
$ws.Range("J110").Value = 2352.7273
$ws.Range("K110").Value = 763.375
$ws.Range("L110").Value = 2352.7273
$ws.Range("M110").Value = 1281.625
$ws.Range("N110").Value = -6442.7273
# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 40552.652
$ws.Range("I116").Value = 65125.312
$ws.Range("K116").Value = 65125.312
$ws.Range("M116").Value = -62831.312
# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3369.6843
$ws.Range("I122").Value = 1823.1428
$ws.Range("J122").Value = 7700
$ws.Range("K122").Value = 5469.428400000001
$ws.Range("L122").Value = 23100
$ws.Range("M122").Value = -3019.428400000001
$ws.Range("N122").Value = -28000

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 40552.652
$ws.Range("I3").Value = 65125.312
$ws.Range("K3").Value = 65125.312
$ws.Range("M3").Value = -65011.312
# row 44 (Leve Item ID 1643)
$ws.Range("H44").Value = 50000
$ws.Range("I44").Value = 50000
$ws.Range("K44").Value = 50000
$ws.Range("M44").Value = -49503
# row 80 (Leve Item ID 13747)
$ws.Range("H80").Value = 822.7368
$ws.Range("I80").Value = 179.8
$ws.Range("J80").Value = 1052.3572
$ws.Range("K80").Value = 179.8
$ws.Range("L80").Value = 1052.3572
$ws.Range("M80").Value = 818.2
$ws.Range("N80").Value = -3048.3572
# row 83 (Leve Item ID 13747)
$ws.Range("H83").Value = 822.7368
$ws.Range("I83").Value = 179.8
$ws.Range("J83").Value = 1052.3572
$ws.Range("K83").Value = 899
$ws.Range("L83").Value = 5261.786
$ws.Range("M83").Value = 4093
$ws.Range("N83").Value = -15245.786
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 11626.637
$ws.Range("I86").Value = 4146.5
$ws.Range("J86").Value = 15901
$ws.Range("K86").Value = 4146.5
$ws.Range("L86").Value = 15901
$ws.Range("M86").Value = -3023.5
$ws.Range("N86").Value = -18147
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 11626.637
$ws.Range("I89").Value = 4146.5
$ws.Range("J89").Value = 15901
$ws.Range("K89").Value = 20732.5
$ws.Range("L89").Value = 79505
$ws.Range("M89").Value = -15116.5
$ws.Range("N89").Value = -90737

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1546.3846
$ws.Range("I16").Value = 1426.3636
$ws.Range("K16").Value = 1426.3636
$ws.Range("M16").Value = -1139.3636
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1831
$ws.Range("I31").Value = 1093.7587
$ws.Range("K31").Value = 1093.7587
$ws.Range("M31").Value = -798.7587000000001
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1831
$ws.Range("I34").Value = 1093.7587
$ws.Range("K34").Value = 1093.7587
$ws.Range("M34").Value = -891.7587000000001
# row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1546.3846
$ws.Range("I113").Value = 1426.3636
$ws.Range("K113").Value = 1426.3636
$ws.Range("M113").Value = 743.6364000000001
# row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2299.36
$ws.Range("I132").Value = 1965.1945
$ws.Range("J132").Value = 3158.6428
$ws.Range("K132").Value = 5895.583500000001
$ws.Range("L132").Value = 9475.928400000001
$ws.Range("M132").Value = -3365.583500000001
$ws.Range("N132").Value = -14535.9284

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 25 (Leve Item ID 4709)
$ws.Range("H25").Value = 475
$ws.Range("I25").Value = 475
$ws.Range("K25").Value = 1425
$ws.Range("M25").Value = -1256
# row 29 (Leve Item ID 4698)
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
# row 30 (Leve Item ID 4709)
$ws.Range("H30").Value = 475
$ws.Range("I30").Value = 475
$ws.Range("K30").Value = 1425
$ws.Range("M30").Value = -1323
# row 75 (Leve Item ID 12863)
$ws.Range("H75").Value = 1992
$ws.Range("J75").Value = 2240
$ws.Range("L75").Value = 6720
$ws.Range("N75").Value = -8716
# row 78 (Leve Item ID 12863)
$ws.Range("H78").Value = 1992
$ws.Range("J78").Value = 2240
$ws.Range("L78").Value = 20160
$ws.Range("N78").Value = -30144

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 1642.8572
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 1416.6666
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 1416.6666
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -3412.6666
# row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 1642.8572
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 1416.6666
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 7083.333000000001
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -17067.333
# row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2130.6428
$ws.Range("I126").Value = 1607
$ws.Range("J126").Value = 2421.5557
$ws.Range("K126").Value = 4821
$ws.Range("L126").Value = 7264.6671
$ws.Range("M126").Value = -2351
$ws.Range("N126").Value = -12204.6671

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 45 (Leve Item ID 3851)
$ws.Range("H45").Value = 5911.1113
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1593
# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3367.2173
$ws.Range("I122").Value = 2563
$ws.Range("J122").Value = 3884.2144
$ws.Range("K122").Value = 7689
$ws.Range("L122").Value = 11652.6432
$ws.Range("M122").Value = -5239
$ws.Range("N122").Value = -16552.6432
# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3327.1924
$ws.Range("I132").Value = 1937.0555
$ws.Range("J132").Value = 6455
$ws.Range("K132").Value = 5811.166499999999
$ws.Range("L132").Value = 19365
$ws.Range("M132").Value = -3281.166499999999
$ws.Range("N132").Value = -24425
# row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 51900
$ws.Range("J133").Value = 51900
$ws.Range("L133").Value = 51900
$ws.Range("N133").Value = -56960
# row 135 (Leve Item ID 42036)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 112473.11
$ws.Range("I122").Value = 126157.25
$ws.Range("K122").Value = 378471.75
$ws.Range("M122").Value = -376021.75

Write-Host "Updated market-data columns across 8 sheets."